# Moderhinke Dashboard update
# - decrements a handful of existing "D" values by 2 (weeks 202440..202504,
#   rows for farms_total_count / farms_to_examine_count)
# - appends a new week block (202505, 2025-02-02) for rows 87-91

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Adjust existing values (each decreases by 2) ---
$ws.Range("D2").Value = 11765
$ws.Range("D3").Value = 11510
$ws.Range("D7").Value = 11869
$ws.Range("D8").Value = 11157
$ws.Range("D12").Value = 11977
$ws.Range("D13").Value = 10655
$ws.Range("D17").Value = 12021
$ws.Range("D18").Value = 10104
$ws.Range("D22").Value = 12078
$ws.Range("D23").Value = 9674
$ws.Range("D27").Value = 12120
$ws.Range("D28").Value = 9213
$ws.Range("D32").Value = 12160
$ws.Range("D33").Value = 8805
$ws.Range("D37").Value = 12198
$ws.Range("D38").Value = 8281
$ws.Range("D42").Value = 12233
$ws.Range("D43").Value = 7713
$ws.Range("D47").Value = 12258
$ws.Range("D48").Value = 7165
$ws.Range("D52").Value = 12286
$ws.Range("D53").Value = 6523
$ws.Range("D57").Value = 12309
$ws.Range("D58").Value = 5913
$ws.Range("D62").Value = 12320
$ws.Range("D63").Value = 5747
$ws.Range("D67").Value = 12342
$ws.Range("D68").Value = 5563
$ws.Range("D72").Value = 12358
$ws.Range("D73").Value = 5127
$ws.Range("D77").Value = 12381
$ws.Range("D78").Value = 4650
$ws.Range("D82").Value = 12399
$ws.Range("D83").Value = 4132

# --- Append the new week (YearWeekIso 202505, LastDayOfWeek 2025-02-02) ---
# Serial date 45690 == 2025-02-02 (same epoch as the existing B column cells)
$newWeek = 202505
$newDateSerial = 45690

$ws.Range("A87").Value = $newWeek
$ws.Range("B87").Value = $newDateSerial
$ws.Range("C87").Value = "farms_total_count"
$ws.Range("D87").Value = 12426

$ws.Range("A88").Value = $newWeek
$ws.Range("B88").Value = $newDateSerial
$ws.Range("C88").Value = "farms_to_examine_count"
$ws.Range("D88").Value = 3725

$ws.Range("A89").Value = $newWeek
$ws.Range("B89").Value = $newDateSerial
$ws.Range("C89").Value = "farms_examined_count"
$ws.Range("D89").Value = 8701

$ws.Range("A90").Value = $newWeek
$ws.Range("B90").Value = $newDateSerial
$ws.Range("C90").Value = "farms_examined_positive_count"
$ws.Range("D90").Value = 1536

$ws.Range("A91").Value = $newWeek
$ws.Range("B91").Value = $newDateSerial
$ws.Range("C91").Value = "farms_examined_negative_count"
$ws.Range("D91").Value = 7165

# Reuse the same date display format as the rest of column B (numFmtId 14, "m/d/yyyy")
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B87:B91").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Match the author's final selection/view state ---
$ws.Range("D90").Select()
